# Verification_Plan.xlsx update - "Added updated SPI Wrapper project files"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level metadata (absolute path recorded by Excel) ---
# Not exposed via a stable, generic COM property on this engine; the
# author's local path is not meaningful application content, so it is
# intentionally left alone here.

# --- Fill in missing "-" cells in column D for existing rows 3-5 ---
$ws.Range("D3").Value2 = "-"
$ws.Range("D4").Value2 = "-"
$ws.Range("D5").Value2 = "-"

# --- New column A width ---
$ws.Columns.Item(1).ColumnWidth = 20

# --- New rows 6-10: requirement rows for the SPI wrapper project ---
# Cell values are entered in this particular order so the shared-string
# table comes out in the same sequence as the authored workbook.

$ws.Range("A6").Value2 = "reset"
$ws.Range("B6").Value2 = "When reset is asserted , MISO , rx_data , rx_valid will be low"
$ws.Range("C6").Value2 = "Directed at the beginning of simulation then randomized"

$ws.Range("A7").Value2 = "SS_n for all cases"
$ws.Range("E7").Value2 = "Assertions "

$ws.Range("A8").Value2 = "SS_n for  read_data"

$ws.Range("A9").Value2 = "tx_valid for read data"

$ws.Range("B8").Value2 = "when the least significant bit [0:2]of array_rand inside {111} and counter not equal 24 , SS_n will be low and if equal will be high"
$ws.Range("B7").Value2 = "when the least significant bit [0:2]of array_rand inside {000,001,110} and counter not equal 14 , SS_n will be low and if equal will be high"
$ws.Range("B9").Value2 = "when the least significant bit [0:2]of array_rand equal 3'b111 and counter equal 23 , then tx_valid will be asserted"

$ws.Range("A10").Value2 = "array_rand for all cases"
$ws.Range("B10").Value2 = "if SS_n fell ,  the least significant bit [0:2]of array_rand will be inside {000,001,110,111}"

# --- Remaining cells in the new rows, reusing already-existing strings ---
$ws.Range("D6").Value2 = "-"
$ws.Range("E6").Value2 = "Reference Model, Assertions"

$ws.Range("C7").Value2 = "Randomized"
$ws.Range("D7").Value2 = "-"

$ws.Range("C8").Value2 = "Randomized"
$ws.Range("D8").Value2 = "-"
$ws.Range("E8").Value2 = "Assertions "

$ws.Range("C9").Value2 = "Randomized"
$ws.Range("D9").Value2 = "-"
$ws.Range("E9").Value2 = "-"

$ws.Range("C10").Value2 = "Randomized"
$ws.Range("D10").Value2 = "-"
$ws.Range("E10").Value2 = "-"

# --- Cell formatting for the new rows ---
# Create the "vertical top only" style first (cellXfs index 5), then the
# "vertical top + wrap text" style (cellXfs index 6), matching the order
# in which the two new formats were introduced.

# Column C rows 7-10: vertical top only (no wrap)
$ws.Range("C7:C10").VerticalAlignment = -4160

# Column B (wrap text, vertical top) for rows 6-10
$ws.Range("B6:B10").WrapText = $true
$ws.Range("B6:B10").VerticalAlignment = -4160

# Column C: row 6 matches column B styling (wrap + vertical top)
$ws.Range("C6").WrapText = $true
$ws.Range("C6").VerticalAlignment = -4160

# --- Row heights for the new rows ---
$ws.Rows.Item(6).RowHeight = 49.2
$ws.Rows.Item(7).RowHeight = 59.4
$ws.Rows.Item(8).RowHeight = 48
$ws.Rows.Item(9).RowHeight = 42.6
$ws.Rows.Item(10).RowHeight = 42.6

# --- Sheet view: zoom + new selection ---
$excel.ActiveWindow.Zoom = 83
[void]$ws.Range("I7").Select()

Write-Output "edit complete"
